# Rename the 4th and 5th sheets (swap their logical roles), then make the
# 4th sheet ("prodTestSheetName" after rename) the active/selected sheet.

$wb = $excel.ActiveWorkbook

$sheet4 = $wb.Worksheets.Item(4)
$sheet5 = $wb.Worksheets.Item(5)

$sheet5.Name = "xxprodTestSheetName"
$sheet4.Name = "prodTestSheetName"

$sheet4.Select()
$sheet4.Activate()
